$d = $word.ActiveDocument

# 1) "May 7, 2018" -> "7 de Maio de 2018"  (inside a hyperlink field, sole run in its paragraph)
$d.Content.Find.Execute("May 7, 2018", $true, $false, $false, $false, $false,
                         $true, 1, $false, "7 de Maio de 2018", 2)

# 2) Title: "Replacing Product Visions with Customer Journey Visions" -> pt-BR
#    (sole run inside its hyperlink/paragraph, Find/Replace preserves formatting fine here)
$d.Content.Find.Execute("Replacing Product Visions with Customer Journey Visions", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Substituindo Visões de Produto por Visões da Jornada do Cliente", 2)

# 3) "Homer Simpson car vision: " -> "Visão do carro de Homer Simpson: "
#    (first run of its paragraph, Find/Replace preserves formatting fine here)
$d.Content.Find.Execute("Homer Simpson car vision: ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Visão do carro de Homer Simpson: ", 2)

# 4) "The Homer" -> "O Homer"
#    This run is the SECOND run in its paragraph (preceded by the run above) and sits
#    inside a w:hyperlink. Using Find.Execute's replace on this particular run loses the
#    run's own character formatting (color / underline) - the engine re-creates the run
#    using the preceding run's rPr instead of its own. Work around it by locating the
#    exact sub-range "by hand" (skip the hyperlink's leading zero-width positions),
#    snapshotting its Font formatting, replacing the text via Range.Text, and then
#    re-applying the snapshotted formatting to the freshly written range.
$targetOld = "The Homer"
$targetNew = "O Homer"

$para = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Contains($targetOld)) {
        $para = $p
        break
    }
}

$hyperlink = $para.Range.Hyperlinks.Item(1)
$hr = $hyperlink.Range.Duplicate

# Find where the visible text actually begins inside the hyperlink's range (the range
# includes some leading zero-width/invisible positions that aren't part of the text).
$visibleStart = -1
for ($s = $hr.Start; $s -lt $hr.End; $s++) {
    $probe = $d.Range($s, $s + 1)
    if ($probe.Text.Length -gt 0) {
        $visibleStart = $s
        break
    }
}

$target = $d.Range($visibleStart, $visibleStart + $targetOld.Length)
$savedColor = $target.Font.Color
$savedUnderline = $target.Font.Underline

$target.Text = $targetNew

$target2 = $d.Range($visibleStart, $visibleStart + $targetNew.Length)
$target2.Font.Color = $savedColor
$target2.Font.Underline = $savedUnderline

# 5) "Where do you see the product in X years?" -> "Onde você vê o produto em X anos?"
#    (sole run of its paragraph, Find/Replace preserves formatting fine here)
$d.Content.Find.Execute("Where do you see the product in X years?", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Onde você vê o produto em X anos?", 2)
